$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new row (96) with the next quarterly data point (01-07-2021).
# Force column A to remain plain text (matching the existing series labels)
# instead of letting Excel auto-convert the date-like string to a date value.
$ws.Range("A96").NumberFormat = "@"
$ws.Range("A96").Value = "01-07-2021"
$ws.Range("A96").Style = "Normal"
$ws.Range("B96").Value = 320
$ws.Range("C96").Value = 141
$ws.Range("D96").Value = 390
$ws.Range("E96").Value = 201
$ws.Range("F96").Value = 313
